$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (so purely-numeric
# looking strings like "130" or "969" stay text, matching the source data),
# then reset the style back to Normal so no stray number-format style sticks
# around on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("A2").Value = "130, 130, 426, 786"
$ws.Range("B2").Value = 1
Set-TextValue $ws.Range("C2") "130"
Set-TextValue $ws.Range("D2") "969"

# Row 3
$ws.Range("A3").Value = "98, 458, 754, 1082"
$ws.Range("B3").Value = 1
Set-TextValue $ws.Range("C3") "1082"
Set-TextValue $ws.Range("D3") "946"

# Row 4
$ws.Range("A4").Value = "98, 130, 786, SF"
$ws.Range("B4").Value = 1
Set-TextValue $ws.Range("C4") "786"
Set-TextValue $ws.Range("D4") "959"

# The remaining combinations (previously rows 5-14) were consolidated into
# the "ionic up"/"ionic down" force rows above, so clear the now-unused
# rows; this shrinks the sheet's used range back down to A1:D4.
$ws.Range("A5:D14").ClearContents()
